$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 106
$ws.Range("H106").Value = 18484.55
$ws.Range("I106").Value = 18484.55
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 18484.55
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -17853.55
$ws.Range("N106").Value = $null
# Row 127
$ws.Range("H127").Value = 869.4211
$ws.Range("I127").Value = 434
$ws.Range("J127").Value = 1046.8148
$ws.Range("K127").Value = 1302
$ws.Range("L127").Value = 3140.4444
$ws.Range("M127").Value = 3658
$ws.Range("N127").Value = -13060.4444
# Row 129
$ws.Range("H129").Value = 955.34485
$ws.Range("I129").Value = 997
$ws.Range("J129").Value = 953.8570999999999
$ws.Range("K129").Value = 2991
$ws.Range("L129").Value = 2861.5713
$ws.Range("M129").Value = 2009
$ws.Range("N129").Value = -12861.5713
# Row 135
$ws.Range("H135").Value = 22728262
$ws.Range("I135").Value = 646.17645
$ws.Range("J135").Value = 100002160
$ws.Range("K135").Value = 5815.58805
$ws.Range("L135").Value = 900019440
$ws.Range("M135").Value = -3280.58805
$ws.Range("N135").Value = -900024510

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 1491.5555
$ws.Range("I110").Value = 1588.5
$ws.Range("J110").Value = 1297.6666
$ws.Range("K110").Value = 1588.5
$ws.Range("L110").Value = 1297.6666
$ws.Range("M110").Value = 456.5
$ws.Range("N110").Value = -5387.6666

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 829.7778
$ws.Range("I64").Value = 600
$ws.Range("J64").Value = 858.5
$ws.Range("K64").Value = 600
$ws.Range("L64").Value = 858.5
$ws.Range("M64").Value = -375
$ws.Range("N64").Value = -1308.5
# Row 67
$ws.Range("H67").Value = 829.7778
$ws.Range("I67").Value = 600
$ws.Range("J67").Value = 858.5
$ws.Range("K67").Value = 600
$ws.Range("L67").Value = 858.5
$ws.Range("M67").Value = 180
$ws.Range("N67").Value = -2418.5
# Row 105
$ws.Range("H105").Value = 885568.75
$ws.Range("I105").Value = 1327753.1
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 1327753.1
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = -1326006.1
$ws.Range("N105").Value = -4694
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").Value = $null

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2884
$ws.Range("I16").Value = 4094.4
$ws.Range("J16").Value = 866.6667
$ws.Range("K16").Value = 4094.4
$ws.Range("L16").Value = 866.6667
$ws.Range("M16").Value = -3807.4
$ws.Range("N16").Value = -1440.6667
# Row 31
$ws.Range("H31").Value = 15969.886
$ws.Range("I31").Value = 20229.076
$ws.Range("J31").Value = 13453.091
$ws.Range("K31").Value = 20229.076
$ws.Range("L31").Value = 13453.091
$ws.Range("M31").Value = -19934.076
$ws.Range("N31").Value = -14043.091
# Row 34
$ws.Range("H34").Value = 15969.886
$ws.Range("I34").Value = 20229.076
$ws.Range("J34").Value = 13453.091
$ws.Range("K34").Value = 20229.076
$ws.Range("L34").Value = 13453.091
$ws.Range("M34").Value = -20027.076
$ws.Range("N34").Value = -13857.091
# Row 105
$ws.Range("H105").Value = 771.9167
$ws.Range("I105").Value = 471.66666
$ws.Range("J105").Value = 1072.1666
$ws.Range("K105").Value = 471.66666
$ws.Range("L105").Value = 1072.1666
$ws.Range("M105").Value = 1275.33334
$ws.Range("N105").Value = -4566.1666
# Row 113
$ws.Range("H113").Value = 2884
$ws.Range("I113").Value = 4094.4
$ws.Range("J113").Value = 866.6667
$ws.Range("K113").Value = 4094.4
$ws.Range("L113").Value = 866.6667
$ws.Range("M113").Value = -1924.4
$ws.Range("N113").Value = -5206.6667

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 222
$ws.Range("I7").Value = 190
$ws.Range("K7").Value = 570
$ws.Range("M7").Value = -458
# Row 80
$ws.Range("H80").Value = 2646.2
$ws.Range("J80").Value = 2761
$ws.Range("L80").Value = 8283
$ws.Range("N80").Value = -10155
# Row 83
$ws.Range("H83").Value = 2646.2
$ws.Range("J83").Value = 2761
$ws.Range("L83").Value = 24849
$ws.Range("N83").Value = -34209
# Row 92
$ws.Range("H92").Value = 331.2
$ws.Range("J92").Value = 331.2
$ws.Range("L92").Value = 993.5999999999999
$ws.Range("N92").Value = -3489.6
# Row 109
$ws.Range("H109").Value = 946.381
$ws.Range("I109").Value = 573.7646999999999
$ws.Range("J109").Value = 2530
$ws.Range("K109").Value = 1721.2941
$ws.Range("L109").Value = 7590
$ws.Range("M109").Value = -681.2940999999998
$ws.Range("N109").Value = -9670
# Row 139
$ws.Range("H139").Value = 11113695
$ws.Range("I139").Value = 2000.4736
$ws.Range("J139").Value = 30306622
$ws.Range("K139").Value = 6001.4208
$ws.Range("L139").Value = 90919866
$ws.Range("M139").Value = -861.4207999999999
$ws.Range("N139").Value = -90930146

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2376.8667
$ws.Range("I61").Value = 2237.111
$ws.Range("K61").Value = 2237.111
$ws.Range("M61").Value = -2035.111
# Row 107
$ws.Range("H107").Value = 2600
$ws.Range("I107").Value = 2600
$ws.Range("K107").Value = 2600
$ws.Range("M107").Value = -680
# Row 113
$ws.Range("H113").Value = 2376.8667
$ws.Range("I113").Value = 2237.111
$ws.Range("K113").Value = 2237.111
$ws.Range("M113").Value = -67.11099999999988
# Row 122
$ws.Range("H122").Value = 2199.5938
$ws.Range("I122").Value = 2122.4783
$ws.Range("J122").Value = 2396.6667
$ws.Range("K122").Value = 6367.4349
$ws.Range("L122").Value = 7190.000100000001
$ws.Range("M122").Value = -3917.4349
$ws.Range("N122").Value = -12090.0001
# Row 136
$ws.Range("H136").Value = 401268.62
$ws.Range("I136").Value = 500977.8
$ws.Range("J136").Value = 2432
$ws.Range("K136").Value = 1502933.4
$ws.Range("L136").Value = 7296
$ws.Range("M136").Value = -1500383.4
$ws.Range("N136").Value = -12396

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = $null
# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = $null
# Row 113
$ws.Range("H113").Value = 1635.75
$ws.Range("I113").Value = 414.33334
$ws.Range("J113").Value = 5300
$ws.Range("K113").Value = 1243.00002
$ws.Range("L113").Value = 15900
$ws.Range("M113").Value = 926.9999800000001
$ws.Range("N113").Value = -20240
